$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.562.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.623.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "197.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.647"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000303"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.196.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "601.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.637.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.623.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("E26").Value = "  -3.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.87%  "

$ws.Range("E31").Value = "  +1.52%  "

$ws.Range("E32").Value = "  -2.79%  "

$ws.Range("E33").Value = "  +2.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0885"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.897.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "538.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.91%  "

$ws.Range("E41").Value = "  -1.07%  "

$ws.Range("E42").Value = "  -2.64%  "

$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0459"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.91%  "

$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("E51").Value = "  +2.12%  "
